# Auto-generated edit script: refresh crypto price/volume data
# (values scraped 2024-08-13 per commit message), including three
# coin-rank swaps that shifted rows 16/17/18, 28/29 and 44/45.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.307.70'
$ws.Range("E2").Value = '  +2.09%  '

$ws.Range("D3").Value = '2.688.06'
$ws.Range("E3").Value = '  -0.78%  '

$ws.Range("E4").Value = '  +0.28%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '524.51'
$ws.Range("E5").Value = '  +1.49%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.10'
$ws.Range("E6").Value = '  +1.50%  '

$ws.Range("E7").Value = '  +0.19%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.577'
$ws.Range("E8").Value = '  +1.69%  '

$ws.Range("D9").Value = '2.710.98'
$ws.Range("E9").Value = '  +0.15%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.49'
$ws.Range("E10").Value = '  +4.13%  '

$ws.Range("E11").Value = '  -0.49%  '

$ws.Range("E12").Value = '  +1.09%  '

$ws.Range("E13").Value = '  +1.60%  '

$ws.Range("D14").Value = '3.168.86'
$ws.Range("E14").Value = '  -0.24%  '

$ws.Range("D15").Value = '60.293.66'
$ws.Range("E15").Value = '  +2.02%  '

# Row 16: coin identity changed
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.37'
$ws.Range("E16").Value = '  +1.49%  '

# Row 17: coin identity changed
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000139'
$ws.Range("E17").Value = '  +0.67%  '

# Row 18: coin identity changed
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '2.708.40'
$ws.Range("E18").Value = '  +0.65%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '351.99'
$ws.Range("E19").Value = '  +1.66%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.54'
$ws.Range("E20").Value = '  -0.20%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.57'
$ws.Range("E21").Value = '  +1.02%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.35'
$ws.Range("E22").Value = '  +2.62%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.07%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.07'
$ws.Range("E24").Value = '  +3.28%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.422'
$ws.Range("E25").Value = '  +0.24%  '

$ws.Range("E26").Value = '  +4.95%  '

$ws.Range("E27").Value = '  +0.34%  '

# Row 28: coin identity changed
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.36'
$ws.Range("E28").Value = '  +1.83%  '

# Row 29: coin identity changed
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").Value = '0.0₃0819'
$ws.Range("E29").Value = '  +0.18%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.86'
$ws.Range("E30").Value = '  +6.41%  '

$ws.Range("E31").Value = '  +0.19%  '

$ws.Range("E32").Value = '  +1.28%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.15'
$ws.Range("E33").Value = '  +0.42%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '147.24'
$ws.Range("E34").Value = '  -1.94%  '

$ws.Range("E35").Value = '  +6.17%  '

$ws.Range("E36").Value = '  +9.20%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.959'
$ws.Range("E37").Value = '  -5.93%  '

$ws.Range("E38").Value = '  +8.40%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.878'
$ws.Range("E39").Value = '  +3.49%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.90'
$ws.Range("E40").Value = '  +0.06%  '

$ws.Range("E41").Value = '  +0.27%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '285.42'
$ws.Range("E42").Value = '  +1.45%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.13'
$ws.Range("E43").Value = '  +0.78%  '

# Row 44: coin identity changed
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0993'
$ws.Range("E44").Value = '  +1.06%  '

# Row 45: coin identity changed
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.614'
$ws.Range("E45").Value = '  -1.35%  '

$ws.Range("E46").Value = '  -0.11%  '

$ws.Range("D47").Value = '2.132.87'
$ws.Range("E47").Value = '  +5.91%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0543'
$ws.Range("E48").Value = '  +1.39%  '

$ws.Range("E49").Value = '  +3.25%  '

$ws.Range("E50").Value = '  +2.16%  '

$ws.Range("E51").Value = '  +1.76%  '

